# Fix up OCR-extracted metadata values on the "metadata" sheet.
# (see commit: "add function to export ocr output to eagle eye excel formate")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("metadata")

# account_number: T6-1806128-3 -> T6-18060128-3
$ws.Range("B2").Value = "T6-18060128-3`n`f"

# account_holder: NUR IZZAHTI BINTI AZEMAN -> NUR IZZAHTI BINT! AZEMAN
$ws.Range("B3").Value = "NUR IZZAHTI BINT! AZEMAN`n`f"

# address: "53100 SELANGOR SELANGOR" -> "53100 SELANGOR, SELANGOR"
$ws.Range("B4").Value = "NO 19 JALAN Nd FASA DA`nTAMAN MELAWATI`n`nKUALA LUMPUR`n`n53100 SELANGOR, SELANGOR`n`f"

# date: "150H2020" -> "{0B 2020"
$ws.Range("B5").Value = "{0B 2020`n`f"

# address (second occurrence): "3419 JALAN Ald FASA OA" -> "7419 JALAN Ald FASA GA"
$ws.Range("B6").Value = "UR IZZAHTI BINT! AZEM`n7419 JALAN Ald FASA GA`n`f"
